# stagingMSRPs.xlsx update: 2021 MY refresh for several Lexus rows, plus
# newly added ES250 / RX-performance / Black-Line-SE trims.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. 2020 -> 2021 model-year bump (+ a handful of MSRP changes) for the
#    existing ES-300h/ES-350 block (rows 10-16) and RX block (rows 55-63).
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = 2021

$ws.Range("C11").Value = 2021
$ws.Range("D11").Value = 45100

$ws.Range("C12").Value = 2021
$ws.Range("D12").Value = 48900

$ws.Range("C13").Value = 2021

$ws.Range("C14").Value = 2021
$ws.Range("D14").Value = 47010

$ws.Range("C15").Value = 2021
$ws.Range("D15").Value = 50810

$ws.Range("C16").Value = 2021
$ws.Range("D16").Value = 45700

$ws.Range("C55").Value = 2021
$ws.Range("D55").Value = 45070

$ws.Range("C56").Value = 2021
$ws.Range("D56").Value = 46470

$ws.Range("C57").Value = 2021
$ws.Range("D57").Value = 47900

$ws.Range("C58").Value = 2021
$ws.Range("D58").Value = 49300

$ws.Range("C59").Value = 2021
$ws.Range("D59").Value = 48550

$ws.Range("C60").Value = 2021
$ws.Range("D60").Value = 49950

$ws.Range("C61").Value = 2021
$ws.Range("D61").Value = 47720

$ws.Range("C62").Value = 2021
$ws.Range("D62").Value = 51110

$ws.Range("C63").Value = 2021
$ws.Range("D63").Value = 51200

# ---------------------------------------------------------------------------
# 2. New rows 81-95: ES 250 trims, RX performance/luxury trims and the new
#    Black Line Special Edition RX trims. Trim names (col B) are entered in
#    the same order the sheet author typed them (82-85,81, then 86-94), and
#    the SE order codes (col A, rows 92/93/94/85) are filled in afterwards.
# ---------------------------------------------------------------------------

$ws.Range("B82").Value = "ES 250 LUXURY"
$ws.Range("B83").Value = "ES 250 ULTRA LUXURY"
$ws.Range("B84").Value = "ES 250 F SPORT"
$ws.Range("B85").Value = "ES 350 F SPORT BLACK LINE SPECIAL EDITION"
$ws.Range("B81").Value = "ES 250"
$ws.Range("B86").Value = "RX 350 F SPORT PERFORMANCE FWD"
$ws.Range("B87").Value = "RX 350 F SPORT PERFORMANCE AWD"
$ws.Range("B88").Value = "RX 450h F-SPORT PERFORMANCE AWD"
$ws.Range("B89").Value = "RX 350L LUXURY FWD"
$ws.Range("B90").Value = "RX 350L LUXURY AWD"
$ws.Range("B91").Value = "RX 450hL LUXURY AWD"
$ws.Range("B92").Value = "RX 350 F SPORT BLACK LINE SPECIAL EDITION"
$ws.Range("B93").Value = "RX 350 F SPORT BLACK LINE SPECIAL EDITION"
$ws.Range("B94").Value = "RX 450h F SPORT AWD BLACK LINE SPECIAL EDITION"

$ws.Range("A92").Value = "9422SE"
$ws.Range("A93").Value = "9426SE"
$ws.Range("A94").Value = "9446SE"
$ws.Range("A85").Value = "9005SE"

# Remaining (non-shared-string) cells: numeric order codes, year, MSRP, DPHF.
$ws.Range("A81").Value = 9012
$ws.Range("A86").Value = 9423
$ws.Range("A87").Value = 9427
$ws.Range("A88").Value = 9447
$ws.Range("A89").Value = 9432
$ws.Range("A90").Value = 9436
$ws.Range("A91").Value = 9457

$newRows = @(
  @{ Row=81; C=2021; D=39900 },
  @{ Row=82; C=2021; D=45100 },
  @{ Row=83; C=2021; D=48900 },
  @{ Row=84; C=2021; D=45700 },
  @{ Row=85; C=2021; D=46550 },
  @{ Row=86; C=2021; D=50950 },
  @{ Row=87; C=2021; D=52350 },
  @{ Row=88; C=2021; D=53520 },
  @{ Row=89; C=2021; D=53900 },
  @{ Row=90; C=2021; D=55300 },
  @{ Row=91; C=2021; D=57110 },
  @{ Row=92; C=2021; D=49235 },
  @{ Row=93; C=2021; D=50635 },
  @{ Row=94; C=2021; D=51885 }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Range("C$row").Value = $r.C
  $ws.Range("D$row").Value = $r.D
  $ws.Range("D$row").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
  $ws.Range("E$row").Value = 1025
  $ws.Range("E$row").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'
}

# Stray styled-but-empty cells trailing the new block (left behind by the
# original author's formatting drag) - columns J/K, rows 91-95.
$ws.Range("J91").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'
$ws.Range("K92").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'
$ws.Range("K93").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'
$ws.Range("K94").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'
$ws.Range("K95").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# ---------------------------------------------------------------------------
# 3. Column B got wider to fit the longer new trim names.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 28.67

# ---------------------------------------------------------------------------
# 4. View state: scrolled down to the newly-added rows, zoomed to 80%, and
#    the active selection moved to A86.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 80
$ws.Range("A86").Select() | Out-Null
